# Atualizações dados 16/07 23h
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 162: match went from "inprogress" (no officials/score yet) to "finished" ---
$ws.Range("G162").Value = "finished"
$ws.Range("H162").Value = 789405
$ws.Range("I162").Value = 784982
$ws.Range("J162").Value = 784886
$ws.Range("N162").Value = 1
$ws.Range("O162").Value = 1

# --- Row 163: match went from "notstarted" to "finished" ---
$ws.Range("G163").Value = "finished"
$ws.Range("H163").Value = 788983
$ws.Range("I163").Value = 784889
$ws.Range("J163").Value = 791416
$ws.Range("N163").Value = 2
$ws.Range("O163").Value = 1

# --- New column P: "dt_insertion" timestamp stamped onto every data row ---
$ws.Range("P1").Value = "dt_insertion"

# Copy the bold/centered/bordered header style from an existing header cell
# onto the new header cell so it matches the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$lastRow = 169
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 16).Value = 45489.90765046296
}

# Establish the datetime number format. Setting the lowercase variant first
# and then the uppercase one mirrors the existing yyyy-mm-dd/YYYY-MM-DD pair
# already present in the workbook's styles (164/165), registering both
# 166 (yyyy-mm-dd h:mm:ss) and 167 (YYYY-MM-DD HH:MM:SS) while only the
# latter ends up applied to cells.
$ws.Range("P2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("P2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("P2:P$lastRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"
